$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (column D) cells. These values must remain stored as
# TEXT (not numbers) to preserve exact formatting (trailing zeros, etc.),
# matching the original workbook where these cells are inline strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.87"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.236"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05862"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.463"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.330"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8091"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8944"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1378"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07222"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03062"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03058"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09325"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.869"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001539"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04702"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006044"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006212"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001263"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004575"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008701"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.559"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.174"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3199"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002341"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03767"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006341"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002500"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007070"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005426"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5404"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.005045"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"

# Update Coin name / Link / Volume(1h) text cells.
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
